$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reagent / source-well / transfer-volume groupings (Source Well -> Reagent, Volume)
$sourceWells = @("A1", "A2", "A3")
$reagents    = @("DNA ligase buffer", "DNA ligase", "BsaI-HF")
$volumes     = @(500, 125, 250)
$destWells   = @("A1", "A2", "A3", "A4")

$uid = 1
$row = 2

for ($s = 0; $s -lt $sourceWells.Length; $s++) {
    for ($d = 0; $d -lt $destWells.Length; $d++) {
        $ws.Cells.Item($row, 1).Value = $uid
        $ws.Cells.Item($row, 2).Value = "level 1 LDV source plate"
        $ws.Cells.Item($row, 3).Value = "384LDV_AQ_B"
        $ws.Cells.Item($row, 4).Value = $sourceWells[$s]
        $ws.Cells.Item($row, 5).Value = "384-Well Level 1 MoClo output plate"
        $ws.Cells.Item($row, 6).Value = "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)"
        $ws.Cells.Item($row, 7).Value = $destWells[$d]
        $ws.Cells.Item($row, 8).Value = $volumes[$s]
        $ws.Cells.Item($row, 9).Value = $reagents[$s]

        $uid++
        $row++
    }
}
